# Insert a new weekly price record for "Vega Modelo de Temuco - Cilantro" at
# row 591. This pushes the existing rows 591-659 down to 592-660 (the
# worksheet grows from 659 to 660 used rows / A1:R659 -> A1:R660).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 591, shifting rows 591:659 down
# to 592:660.
$ws.Rows.Item(591).Insert()

# Populate the newly inserted row with the new data point. Columns that are
# constant across every record in this sheet (Mercado ID, Mercado, Region,
# Codreg, Categoria ID, Categoria, Variedad, Calidad, Unidad de
# comercializacion, Origen, Kg o Unidades, Clasificacion) are carried over
# unchanged; only the date (Fecha) and the measured values (Volumen, Precio
# minimo/maximo/promedio ponderado, Precio $/Kg) differ for this new entry.
$ws.Cells.Item(591, 1).Value = 10
$ws.Cells.Item(591, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(591, 3).Value = "La Araucanía"
$ws.Cells.Item(591, 4).Value = 45212
$ws.Cells.Item(591, 5).Value = 9
$ws.Cells.Item(591, 6).Value = 100112040
$ws.Cells.Item(591, 7).Value = "Cilantro"
$ws.Cells.Item(591, 8).Value = "Sin especificar"
$ws.Cells.Item(591, 9).Value = "Primera"
$ws.Cells.Item(591, 10).Value = 70
$ws.Cells.Item(591, 11).Value = 4000
$ws.Cells.Item(591, 12).Value = 4000
$ws.Cells.Item(591, 13).Value = 4000
$ws.Cells.Item(591, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(591, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(591, 16).Value = 2000
$ws.Cells.Item(591, 17).Value = 2
$ws.Cells.Item(591, 18).Value = "Hortaliza"
